$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header label for column O
$ws.Range("O2").Value = "Utility (Percent)"

# Append " msec" to the latency values in columns I, J, K for data rows 3-38
for ($row = 3; $row -le 38; $row++) {
    foreach ($col in @("I", "J", "K")) {
        $cell = $ws.Range("$col$row")
        $current = $cell.Value2
        $cell.Value = "$current msec"
    }
}
